$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content updates (column B written before column A so the
#     shared-string table is built in the same order as the target file) ---
$ws.Range("B2").Value = "Your Objective"
$ws.Range("A2").Value = "Heder"

$ws.Range("B3").Value = "Achieve your goal"
$ws.Range("A3").Value = "Heder"

$ws.Range("B4").Value = "Choose your pack"
$ws.Range("A4").Value = "Heder"

$ws.Range("B5").Value = "Pricing"
$ws.Range("A5").Value = "Heder"

$ws.Range("B6").Value = "IT-Platforma"
$ws.Range("A6").Value = "Heder"

$ws.Range("B7").Value = "Robert'); DROP TABLE Students;--"
$ws.Range("A7").Value = "XSS"

$ws.Range("B8").Value = "Nice site,  I think I'll take it."
$ws.Range("A8").Value = "XSS"

# --- Column widths: col A narrower (target ~10.89 chars), col B wide
#     enough to show the new, longer strings (target ~57.11 chars).
#     The runtime quantizes ColumnWidth to whole pixels (~1/6 character
#     steps), so these inputs are chosen to land on the pixel bucket
#     closest to the target stored width. ---
$ws.Columns("A").ColumnWidth = 10
$ws.Columns("B").ColumnWidth = 56.3

# --- Selection moved to B10 ---
$ws.Range("B10").Select() | Out-Null

# --- Page setup: portrait, paper size 9 (A4) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
